$d = $word.ActiveDocument

# Update the date line (unique text, safe to use Find & Replace)
$d.Content.Find.Execute("2025-08-30 Saturday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-08-31 Sunday", 2)

# Update the division-problem table cells directly by (row, column) so that
# duplicate cell values (e.g. "855÷9=95, 0" appearing twice) are replaced
# independently with their own new values.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "477÷5=95, 2" },
    @{ Row = 1;  Col = 2; New = "869÷6=144, 5" },
    @{ Row = 1;  Col = 3; New = "837÷2=418, 1" },
    @{ Row = 1;  Col = 4; New = "935÷2=467, 1" },
    @{ Row = 1;  Col = 5; New = "535÷4=133, 3" },

    @{ Row = 5;  Col = 1; New = "702÷8=87, 6" },
    @{ Row = 5;  Col = 2; New = "505÷5=101, 0" },
    @{ Row = 5;  Col = 3; New = "777÷6=129, 3" },
    @{ Row = 5;  Col = 4; New = "420÷9=46, 6" },
    @{ Row = 5;  Col = 5; New = "183÷4=45, 3" },

    @{ Row = 9;  Col = 1; New = "661÷6=110, 1" },
    @{ Row = 9;  Col = 2; New = "453÷6=75, 3" },
    @{ Row = 9;  Col = 3; New = "635÷9=70, 5" },
    @{ Row = 9;  Col = 4; New = "672÷7=96, 0" },
    @{ Row = 9;  Col = 5; New = "489÷8=61, 1" },

    @{ Row = 13; Col = 1; New = "730÷5=146, 0" },
    @{ Row = 13; Col = 2; New = "916÷9=101, 7" },
    @{ Row = 13; Col = 3; New = "868÷8=108, 4" },
    @{ Row = 13; Col = 4; New = "287÷6=47, 5" },
    @{ Row = 13; Col = 5; New = "152÷3=50, 2" },

    @{ Row = 17; Col = 1; New = "693÷9=77, 0" },
    @{ Row = 17; Col = 2; New = "514÷8=64, 2" },
    @{ Row = 17; Col = 3; New = "373÷7=53, 2" },
    @{ Row = 17; Col = 4; New = "519÷2=259, 1" },
    @{ Row = 17; Col = 5; New = "811÷4=202, 3" }
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}
